# Add the four new Specialist Mission rows (580-583) to the bottom of the
# "Specialist_Mission_BoBs.csv" sheet, then move the active selection down
# to reflect the new scroll position, matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mission 580: Elyon / Wowbadger
$ws.Range("A308").Value = 580
$ws.Range("B308").Value = "Elyon"
$ws.Range("C308").Value = "Wowbadger"

# Mission 581: Zephyr / Bastien
$ws.Range("A309").Value = 581
$ws.Range("B309").Value = "Zephyr"
$ws.Range("C309").Value = "Bastien"

# Mission 582: Photurius / Abh
$ws.Range("A310").Value = 582
$ws.Range("B310").Value = "Photurius"
$ws.Range("C310").Value = "Abh"

# Mission 583: Banshee / JUB (new player name, appended to shared strings)
$ws.Range("A311").Value = 583
$ws.Range("B311").Value = "Banshee"
$ws.Range("C311").Value = "JUB"

# Reflect the author's final cursor position after adding the rows.
$ws.Range("E295").Select()
